$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '24.954.45'
$ws.Range("E2").Value = '  +2.11%  '
$ws.Range("D3").Value = '1.701.27'
$ws.Range("E3").Value = '  +0.77%  '
$ws.Range("E4").Value = '  +0.19%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '315.81'
$ws.Range("E5").Value = '  +0.10%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '1.002'
$ws.Range("E6").Value = '  +0.23%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.3974'
$ws.Range("E7").Value = '  +1.81%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.4026'
$ws.Range("E8").Value = '  -0.27%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '1.466'
$ws.Range("E9").Value = '  -1.54%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '52.76'
$ws.Range("E10").Value = '  +0.53%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '1.001'
$ws.Range("E11").Value = '  +0.17%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.08801'
$ws.Range("E12").Value = '  +0.28%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '26.14'
$ws.Range("E13").Value = '  -1.65%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '7.450'
$ws.Range("E14").Value = '  -0.40%  '
$ws.Range("B15").Value = 'ShibaInu'
$ws.Range("C15").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.00001350'
$ws.Range("E15").Value = '  +0.17%  '
$ws.Range("B16").Value = 'Chainlink'
$ws.Range("C16").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '7.958'
$ws.Range("E16").Value = '  -2.82%  '
$ws.Range("D17").Value = '1.708.04'
$ws.Range("E17").Value = '  +1.16%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '96.26'
$ws.Range("E18").Value = '  -1.86%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.07185'
$ws.Range("E19").Value = '  -0.87%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '20.57'
$ws.Range("E20").Value = '  +0.98%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '7.320'
$ws.Range("E21").Value = '  +0.32%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '1.002'
$ws.Range("E22").Value = '  +0.24%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '14.35'
$ws.Range("E23").Value = '  +0.74%  '
$ws.Range("D24").Value = '24.962.50'
$ws.Range("E24").Value = '  +2.10%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.352'
$ws.Range("E25").Value = '  +0.72%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.929'
$ws.Range("E26").Value = '  -3.75%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '23.76'
$ws.Range("E27").Value = '  +4.95%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '6.137'
$ws.Range("E28").Value = '  +14.83%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '161.86'
$ws.Range("E29").Value = '  -3.33%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '150.09'
$ws.Range("E30").Value = '  +8.21%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '8.338'
$ws.Range("E31").Value = '  -3.02%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '2.651'
$ws.Range("E32").Value = '  +25.62%  '
$ws.Range("D33").Value = '1.897.52'
$ws.Range("E33").Value = '  +1.18%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.08540'
$ws.Range("E34").Value = '  -2.75%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.03143'
$ws.Range("E35").Value = '  +3.78%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.037'
$ws.Range("E36").Value = '  -1.33%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '7.128'
$ws.Range("E37").Value = '  -2.42%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.2858'
$ws.Range("E38").Value = '  +2.53%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '10.91'
$ws.Range("E39").Value = '  +0.41%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.09556'
$ws.Range("E40").Value = '  +4.64%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.8216'
$ws.Range("E41").Value = '  +2.37%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '13.98'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.482'
$ws.Range("E43").Value = '  +0.67%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '17.12'
$ws.Range("E44").Value = '  -3.38%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '2.675'
$ws.Range("E45").Value = '  +0.85%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.7378'
$ws.Range("E46").Value = '  +1.72%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '4.256'
$ws.Range("E47").Value = '  -0.20%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.400'
$ws.Range("E48").Value = '  -1.98%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.08765'
$ws.Range("E49").Value = '  +8.60%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.002'
$ws.Range("E50").Value = '  +0.27%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '139.26'
$ws.Range("E51").Value = '  -0.09%  '
